$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at position 13 (shifts 13..23 down to 14..24) ---
$ws.Rows.Item(13).Insert()

# --- Step 2: new row 13 should only have B13/C13 (professor name), no A13 ---
# Copy the B/C formatting from row 14 (which already carries the correct wrap styles)
# down onto the blank row 13, then set values and clear the stray A13 cell that
# Insert() carried over from the old row 13.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("A13").Clear()

# --- Step 3: fix the text that is now misaligned / needs updating ---

# Objetivos: (row 10) -> new Portuguese objectives text
$objPt = 'Apresentar aos alunos o panorama geral da administração estratégica de marketing, capacitando-os a atuar no processo gerencial de marketing sob as perspectivas estratégica e operacional.'
$ws.Range("B10").Value = $objPt
$ws.Range("C10").Value = $objPt

# Programa resumido: (row 14) -> new short-syllabus (PT) text
$resumidoPt = 'Marketing estratégico e planejamento estratégico orientado para o mercado.'
$ws.Range("B14").Value = $resumidoPt
$ws.Range("C14").Value = $resumidoPt

# Programa: (row 16) -> new full syllabus (PT) text
$programaPt = 'MARKETING ESTRATÉGICO1. Evolução do Conceito de Marketing e Sistema de Marketing2. Marketing, Conceito de Valor, Orientação para Mercado3. Análise de Mercado e Comportamento do Consumidor4. Sistema de Informações de Marketing e Inteligência de Mercado5. Modalidades de Marketing6. Fundamentos de Estratégia Empresarial e Marketing Estratégico7. Administração Estratégica e Marketing Estratégico8. Instrumentos Analíticos para Avaliar Oportunidades de Mercado9. Segmentação de Mercado e Posicionamento10. O Plano Estratégico de Marketing11. Comunicação: assessorias de imprensa, SAC''S, Ombudsman'
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# Metodo: (row 19) -> method text
$metodo = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Criterio: (row 20) -> new criterio text
$criterio = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.'
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperacao: (row 21) -> norma text
$norma = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Bibliografia: (row 22) -> new bibliography text
$bib = 'KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM.'
$ws.Range("B22").Value = $bib
$ws.Range("C22").Value = $bib
